# Splits the single-paragraph answer list "0;15;30;60;120;150" into six
# separate numbered list items ("1 = 0 ;", "2 = 15;", ... "6 = 150") for
# both survey questions ("nauka stacjonarna" and "nauka zdalna") that
# contain this text.

$d = $word.ActiveDocument

function Find-ParaIndexByExactText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like "$text*") {
            return $i
        }
    }
    return -1
}

function Split-ListItem($paraIndex, $items) {
    # Re-purpose the existing paragraph (and its pPr/numPr/rPr) for the
    # first item, then insert one new paragraph per remaining item right
    # after it, copying the same paragraph formatting each time.
    $para = $d.Paragraphs($paraIndex)
    $rng = $para.Range
    $innerRng = $d.Range($rng.Start, $rng.End - 1)
    $innerRng.Text = $items[0]

    $prevPara = $d.Paragraphs($paraIndex)
    for ($i = 1; $i -lt $items.Count; $i++) {
        $prevRng = $prevPara.Range
        $prevRng.InsertParagraphAfter()
        $newPara = $d.Paragraphs($paraIndex + $i)
        $newRng = $newPara.Range
        $newInnerRng = $d.Range($newRng.Start, $newRng.End - 1)
        $newInnerRng.Text = $items[$i]
        $prevPara = $newPara
    }
}

$itemsStacjonarna = @("1 = 0 ;", "2 = 15;", "3 = 30;", "4 = 60;", "5 = 120;", "6 = 150")
$itemsZdalna      = @("1 = 0;",  "2 = 15;", "3 = 30;", "4 = 60;", "5 = 120;", "6 = 150")

$idx1 = Find-ParaIndexByExactText("0;15;30;60;120;150")
Split-ListItem $idx1 $itemsStacjonarna

$idx2 = Find-ParaIndexByExactText("0;15;30;60;120;150")
Split-ListItem $idx2 $itemsZdalna

Write-Host "Done. Total paragraphs:" $d.Paragraphs.Count
